# Update "paises.xlsx" COVID country stats + reorder a few country rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 10:04"

# --- Swap country names so alphabetical order is corrected ---
# Kazajistan now sorts before Moldavia
$ws.Range("A59").Value = "Kazajistan"
$ws.Range("A60").Value = "Moldavia"

# Eslovaquia now sorts before Eslovenia
$ws.Range("A88").Value = "Eslovaquia"
$ws.Range("A89").Value = "Eslovenia"

# Seychelles now sorts before Montserrat
$ws.Range("A205").Value = "Seychelles"
$ws.Range("A206").Value = "Montserrat"

# --- Refresh numeric stats (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1292879
$ws.Range("C4").Value = 256
$ws.Range("E4").Value = 998686
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 76942

# Polonia (row 36)
$ws.Range("D36").Value = 5184
$ws.Range("E36").Value = 9108

# Rumania (row 37)
$ws.Range("E37").Value = 7464
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 891

# Ucrania (row 38)
$ws.Range("B38").Value = 14195
$ws.Range("C38").Value = 504
$ws.Range("D38").Value = 2706
$ws.Range("E38").Value = 11128
$ws.Range("G38").Value = 21
$ws.Range("H38").Value = 361

# Row 59 (now Kazajistan)
$ws.Range("B59").Value = 4753
$ws.Range("C59").Value = 175
$ws.Range("D59").Value = 1518
$ws.Range("E59").Value = 3204
$ws.Range("F59").Value = 31
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 31

# Row 60 (now Moldavia)
$ws.Range("B60").Value = 4605
$ws.Range("D60").Value = 1747
$ws.Range("E60").Value = 2713
$ws.Range("F60").Value = 237
$ws.Range("H60").Value = 145

# Row 88 (now Eslovaquia)
$ws.Range("B88").Value = 1455
$ws.Range("C88").Value = 10
$ws.Range("D88").Value = 905
$ws.Range("E88").Value = 524
$ws.Range("F88").Value = 4
$ws.Range("H88").Value = 26

# Row 89 (now Eslovenia)
$ws.Range("B89").Value = 1449
$ws.Range("D89").Value = 247
$ws.Range("E89").Value = 1103
$ws.Range("F89").Value = 13
$ws.Range("H89").Value = 99

# Lituania (row 90)
$ws.Range("B90").Value = 1436
$ws.Range("C90").Value = 3
$ws.Range("D90").Value = 765
$ws.Range("E90").Value = 622

# Letonia (row 96)
$ws.Range("B96").Value = 928
$ws.Range("C96").Value = 19
$ws.Range("E96").Value = 446
$ws.Range("F96").Value = 2

# Taiwan (row 125)
$ws.Range("D125").Value = 355
$ws.Range("E125").Value = 79

# Montenegro (row 132)
$ws.Range("D132").Value = 266
$ws.Range("E132").Value = 50

# Row 205 (now Seychelles)
$ws.Range("D205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("H205").Value = 0

# Row 206 (now Montserrat)
$ws.Range("D206").Value = 7
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1

Write-Host "Edit applied"
